$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Number formats -------------------------------------------------
# Numeric columns (TPS, Energy Use per Transaction, Nakamoto Coefficient)
# use a custom 5-decimal format. Column F/G (Strengths/Weaknesses) use
# text format so "N/A" literals stay left-aligned text. Column E (%) uses
# a whole-number percent format.
$ws.Range("B2:D10").NumberFormat = "0.00000"
$ws.Range("F2:G10").NumberFormat = "@"
$ws.Range("E2:E10").NumberFormat = "0%"

# --- Row 2: Proof of Work -------------------------------------------
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 707
$ws.Range("D2").Value = "N/A"
$ws.Range("E2").Value = 0.25
$ws.Range("F2").Value = "N/A"
$ws.Range("G2").Value = "N/A"

# --- Row 3: Proof of Stake -------------------------------------------
$ws.Range("B3").Value = 30
$ws.Range("C3").Formula = "=AVERAGE(0.03,0.001,0.02,0.04)"
$ws.Range("D3").Value = 379886
$ws.Range("E3").Value = 0.33
$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "N/A"

# --- Row 4: Delegated Proof of Stake ---------------------------------
$ws.Range("B4").Value = "N/A"
$ws.Range("C4").Value = "N/A"
$ws.Range("D4").Value = "N/A"
$ws.Range("E4").Value = "N/A"
$ws.Range("F4").Value = "N/A"
$ws.Range("G4").Value = "N/A"

# --- Row 5: Proof of History ------------------------------------------
$ws.Range("B5").Value = 4501
$ws.Range("C5").Value = 0.001
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = "N/A"
$ws.Range("F5").Value = "N/A"
$ws.Range("G5").Value = "N/A"

# --- Row 6: Proof of Stake with Byzantine Fault Tolerance -------------
$ws.Range("B6").Value = "N/A"
$ws.Range("C6").Value = 0.001
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = 0.33
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "N/A"

# --- Row 7: Proof of History with Proof of Stake ----------------------
$ws.Range("B7").Value = 4501
$ws.Range("C7").Value = "N/A"
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = "N/A"
$ws.Range("F7").Value = "N/A"
$ws.Range("G7").Value = "N/A"

# --- Row 8: zk-proof ---------------------------------------------------
$ws.Range("B8").Value = "N/A"
$ws.Range("C8").Value = "N/A"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "N/A"
$ws.Range("F8").Value = "N/A"
$ws.Range("G8").Value = "N/A"

# --- Row 9: Sharding -----------------------------------------------------
$ws.Range("B9").Value = "N/A"
$ws.Range("C9").Value = "N/A"
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "N/A"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "N/A"

# --- Row 10: DAGs -----------------------------------------------------
$ws.Range("B10").Value = 160000
$ws.Range("C10").Value = "N/A"
$ws.Range("D10").Value = "N/A"
$ws.Range("E10").Value = 0.28
$ws.Range("F10").Value = "N/A"
$ws.Range("G10").Value = "N/A"

# Leave the cursor where the author's session ended up.
$ws.Range("D16").Select() | Out-Null
